$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated simulation results for pl_mw.xlsx (case with 380 kV done)
# Row data: columns B,C,D,F,G,I,K,L,M,N for rows 2-25 (E,H,J,O stay 0)
$colNames = @("B","C","D","F","G","I","K","L","M","N")

$data = @{
    2 = @(0.8981769492673095, 0.05328741999655051, 0.1055607384944341, 1.970185246353921, 0.002542731307204313, 1.626900772744179, 0.6891005233039493, 0.2670590979544798, 0.2483824941078154, 2.940559517220265)
    3 = @(0.861019313806537, 0.04678299856288959, 0.105627226852441, 1.958078412900875, 0.002546479524074944, 1.625926152842403, 0.6472480543886832, 0.2637901141088577, 0.2412460201105446, 2.955287976808911)
    4 = @(0.8387139670163322, 0.04277399513999569, 0.1056669051857178, 1.951566145879241, 0.002548904105026192, 1.625970352547895, 0.6219394267102132, 0.2619130450521112, 0.2370013628938921, 2.965081986762755)
    5 = @(0.8297526878392603, 0.0411363704537564, 0.1056827901799302, 1.949143948640966, 0.002549923208274265, 1.626149952090621, 0.6117237748372872, 0.2611808733482732, 0.2353061666233778, 2.96926182628372)
    6 = @(0.8282724308231764, 0.04086420467980645, 0.1056854108188432, 1.948755730636364, 0.002550094308754535, 1.626189533901972, 0.6100333821906929, 0.2610612757915405, 0.2350267677932436, 2.969967280747369)
    7 = @(0.8385925922116542, 0.04275192555833485, 0.1056671205635755, 1.951532541638656, 0.002548917723316576, 1.625972120438099, 0.6218012587698638, 0.2619030380819325, 0.2369783609897311, 2.965137593591869)
    8 = @(0.885259329862464, 0.05104781000660807, 0.1055839037900395, 1.965819518759659, 0.002543998188788316, 1.626431341061469, 0.6745890561430201, 0.2659049688755388, 0.2458933924613582, 2.945482049076688)
    9 = @(0.9808137101737771, 0.0671994854073148, 0.105411450020922, 2.00115461717904, 0.002535323798408985, 1.632432798296826, 0.7811977434169251, 0.2747842584782063, 0.2644634012544884, 2.91289684671635)
    10 = @(1.053486809619557, 0.07900244529740519, 0.1052788575960548, 2.031593030235044, 0.002529537574253113, 1.639957180516333, 0.8614252583097937, 0.2819367720451424, 0.2787707694379975, 2.89259378335511)
    11 = @(1.087085707047436, 0.08435973594082213, 0.1052172093745902, 2.046416453798997, 0.002527031381026747, 1.644058024262648, 0.898340142443061, 0.2853272864154377, 0.2854240629363929, 2.884148237409121)
    12 = @(1.099886282529667, 0.0863867953660673, 0.1051936698670417, 2.052170379479364, 0.002526100367929503, 1.645708468201207, 0.912379249309339, 0.2866308429793918, 0.2879643018494917, 2.881063894214179)
    13 = @(1.097126009966019, 0.08595030298549489, 0.1051987482278935, 2.050924912251034, 0.002526300077770752, 1.645348676429478, 0.9093529973035572, 0.2863492256137476, 0.2874162923520274, 2.881723099063223)
    14 = @(1.088137267878267, 0.08452653551975686, 0.1052152766840635, 2.046887013691602, 0.002526954425068566, 1.644191852052288, 0.8994939409245148, 0.2854341373025733, 0.2856326336698771, 2.883892205064171)
    15 = @(1.082641478477001, 0.08365422708298809, 0.1052253753966816, 2.044432001175949, 0.002527357576922926, 1.643495968360554, 0.893462829227559, 0.2848761764403349, 0.2845427963260647, 2.885235669700634)
    16 = @(1.051301883872981, 0.07865210093538622, 0.1052828593837969, 2.03064395375597, 0.002529703887445834, 1.639702827584102, 0.8590212141149323, 0.2817179442090207, 0.2783388712879855, 2.893161636818732)
    17 = @(1.032214115718347, 0.07558046280912833, 0.1053177806329657, 2.022435729572692, 0.002531175478074486, 1.637549538491939, 0.837999704372919, 0.279815487823754, 0.274570021605868, 2.898226505015785)
    18 = @(1.021286131145757, 0.07381260926945288, 0.1053377412879062, 2.017806511494584, 0.002532033761938989, 1.636374833454902, 0.8259481251744774, 0.278734123231132, 0.2724159106222004, 2.901214072577091)
    19 = @(1.017594833019473, 0.07321384799539032, 0.1053444782346649, 2.016254923306846, 0.002532326402676182, 1.635988056358364, 0.8218744412001229, 0.2783702044694536, 0.2716889081673486, 2.902238383859341)
    20 = @(1.03424078548602, 0.07590756010813493, 0.1053140761786828, 2.023299993007527, 0.002531017597473867, 1.637772155324541, 0.8402333984373058, 0.2800166749429422, 0.2749698115608439, 2.897679641114607)
    21 = @(1.090775380090349, 0.08494477409693957, 0.1052104271838328, 2.048069224873174, 0.002526761738895185, 1.644528991679536, 0.9023881498882531, 0.2857023880657152, 0.2861559737931003, 2.883251996365814)
    22 = @(1.128175082312737, 0.09084165345410611, 0.1051415505304174, 2.065077044027447, 0.00252408532862686, 1.649513558747174, 0.9433610040334202, 0.2895328076176611, 0.2935878867899362, 2.874486088136365)
    23 = @(1.10817294355553, 0.08769521530071245, 0.1051784162725333, 2.055924598732815, 0.002525504197397034, 1.646801157931606, 0.92146089418776, 0.287477975519522, 0.2896102668322911, 2.879103870587116)
    24 = @(1.03332438532243, 0.07575968543287104, 0.1053157513226317, 2.022908979789108, 0.002531088937109253, 1.637671313229347, 0.8392234395500395, 0.2799256796739513, 0.2747890271133215, 2.897926642447416)
    25 = @(0.9545305116469365, 0.06284174423691979, 0.1054591209063691, 1.99081038805268, 0.002537566948542063, 1.630262601147668, 0.7520244873197441, 0.2722717398388852, 0.2593231880477589, 2.921073591394816)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $colNames.Length; $i++) {
        $cellAddr = "$($colNames[$i])$row"
        $ws.Range($cellAddr).Value = $values[$i]
    }
}
